$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: update A22 formula to use absolute references
$ws.Range("A22").Formula = '=($F$2*$C$8*SIN($E$7-$J$2-$E$8)+$G$2*$C$9*SIN($E$7-$K$2-$E$9))-2*$E$2*$C$7*SIN($I$2)-$C$12*SIN($E$7-$E$12)'

# New cells in row 22
$ws.Range("D22").Formula = '=($F$2*$C$8*SIN($E$7-$J$2-$E$8))'
$ws.Range("E22").Formula = '=$G$2*$C$9*SIN($E$7-$K$2-$E$9)'
$ws.Range("F22").Formula = '=-2*$E$2*$C$7*SIN($I$2)'
$ws.Range("G22").Formula = '=-$C$12*SIN($E$7-$E$12)'

# New cells in row 23
$ws.Range("D23").Formula = '=D22+F22'
$ws.Range("E23").Formula = '=D23+G22'

# Update selection
$ws.Range("M32").Select()
